# Commit as of 4th may2020
# Refresh the FaxAddressBookData test-data scenario with a new run's values:
#   9987288 -> 91827 (and the typo'd 918279 / 991827 variants that appear in
#   the concatenated "recipients" strings), 1000 -> 10384, and the deleted
#   sample recipient "sample2 delete<12345>" -> "Palak Garg<9917186286>".
# A couple of new rows are appended to CreateRecipient / DeleteRecipient for
# a second ("Data2") fax number, and the active sheet/selections move around.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: CreateRecipient
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("CreateRecipient")
$ws1.Range("C2").Value = "'918279"
$ws1.Range("A3").Value = "Data2"
$ws1.Range("C3").Value = "'91827"

# ---------------------------------------------------------------------
# Sheet: EditRecipient
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("EditRecipient")
$ws2.Range("C2").Value = "'91827"
$ws2.Columns.Item(2).ColumnWidth = 11.6

# ---------------------------------------------------------------------
# Sheet: DeleteRecipient
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("DeleteRecipient")
$ws3.Range("C2").Value = "'91827"
$ws3.Range("A5").Value = "Data2"
$ws3.Range("C5").Value = "'91827"
$ws3.Range("D5").Value = "deleted"
$ws3.Columns.Item(4).ColumnWidth = 15.5

# ---------------------------------------------------------------------
# Sheet: AddressCreate
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("AddressCreate")
$ws4.Range("B2").Value = "'10384"
$ws4.Range("C2").Value = "'TrialData Recipient<91827>"
$ws4.Range("D2").Value = "'Palak Garg<9917186286>"
$ws4.Range("E2").Value = "'TrialData Recipient<991827>,Palak Garg<9917186286>"

# ---------------------------------------------------------------------
# Sheet: EditAddressBook
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("EditAddressBook")
$ws5.Range("B2").Value = "'10384"
$ws5.Range("C2").Value = "'TrialData Recipient<91827>,Palak Garg<9917186286>"
$ws5.Range("E2").Value = "'Recepient Updated Recipient<91827>,Palak Garg<9917186286>"
$ws5.Columns.Item(4).ColumnWidth = 27.5

# ---------------------------------------------------------------------
# Sheet: DeleteAddressBook
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("DeleteAddressBook")
$ws6.Range("B2").Value = "'10384"
$ws6.Range("C2").Value = "'Recepient Updated Recipient<91827>,Palak Garg<9917186286>"

# ---------------------------------------------------------------------
# Selections: set each sheet's own cursor position, then finish on
# CreateRecipient so it ends up the active tab (matches the target,
# where CreateRecipient's sheetView gains tabSelected and DeleteRecipient
# loses it).
# ---------------------------------------------------------------------
$ws6.Range("E5").Select()
$ws5.Range("E6").Select()
$ws4.Range("E6").Select()
$ws3.Range("A7").Select()
$ws2.Range("C8").Select()
$ws1.Range("F14").Select()
